$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Change 1 (week 2, row for "Nguyễn Tấn Dũng") ---
# Cell has two paragraphs: "Nguyễn Tấn Dũng" + an empty paragraph.
# Collapse to a single paragraph reading "Nguyễn Ba Duy " (trailing space).
$cell1 = $t.Rows.Item(5).Cells.Item(2)
$cell1.Range.Paragraphs.Item(2).Range.Delete()
$cell1.Range.Paragraphs.Item(1).Range.Text = "Nguyễn Ba Duy "

# --- Change 2 (week 2, row for "Bùi Xuân An") ---
# Cell has "Bùi Xuân An" then an empty second paragraph; add a new run
# "Nguyễn Tấn Dũng" (sz 28) into that empty paragraph.
$cell2 = $t.Rows.Item(7).Cells.Item(2)
$target_start = $cell2.Range.Paragraphs.Item(1).Range.Start
$p = $d.Paragraphs.First
while ($p -ne $null) {
  if ($p.Range.Start -eq $target_start) {
    $p2 = $p.Next()
    $p2.Range.Text = "Nguyễn Tấn Dũng"
    $p2.Range.Font.Size = 14
    break
  }
  $p = $p.Next()
}

# --- Change 3 (week 2, row for "Nguyễn Ba Duy") ---
# Simple text swap: "Nguyễn Ba Duy" -> "Nguyễn Tấn Dũng"
$t.Rows.Item(8).Cells.Item(2).Range.Text = "Nguyễn Tấn Dũng"

# --- Change 4 (week 3, row for "Cả nhóm") ---
# Simple text swap: "Cả nhóm" -> "Nguyễn Tấn Dũng"
$t.Rows.Item(9).Cells.Item(2).Range.Text = "Nguyễn Tấn Dũng"
